# Master_Quest.xlsx / STATUS sheet update
# Row 2: "Alice" -> "Eric", stats reset to 0 / 0%
# Row 3: "John" stats reset to 0 / 0% (name unchanged)
#
# Columns B, E, F hold percentages stored as literal text (e.g. "0%"),
# while C, D hold real numbers. Assigning a string like "0%" directly
# to a General-formatted cell makes Excel auto-convert it into a
# numeric percentage (changing both the stored type and the cell's
# number format/style). To keep the text literal "0%" while preserving
# the original style/border (s="4"), we temporarily force a text
# number format, assign the value, then restore the original
# appearance by pasting the formats from an untouched donor cell in
# the same row (same style) via PasteSpecial.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STATUS")

# --- Row 2: Alice -> Eric -------------------------------------------------
$ws.Range("A2").Value = "Eric"

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "0%"

$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0%"

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "0%"

# --- Row 3: John -----------------------------------------------------------
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "0%"

$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0%"

$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "0%"

# --- Restore original cell formatting/style on the text cells -------------
# C2/C3 were never touched format-wise (only their value changed), so they
# still carry the original style (border, General number format) that B, E
# and F should also end up with.
$ws.Range("C2").Copy()
$ws.Range("B2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("F2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("C3").Copy()
$ws.Range("B3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("F3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$excel.CutCopyMode = $false
